$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the "Meta description: ..." paragraph that currently sits
#    right after the title heading.
# ------------------------------------------------------------------
$metaOld = 'Meta description: Discover the features of Atlantean GigaRise, a highly volatile slot game with up to 294,912 ways to win. Play for free and read our review to learn more.'
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text.StartsWith($metaOld)) {
        $p.Range.Delete()
        break
    }
}

# ------------------------------------------------------------------
# 2) Insert a new bold "Play Atlantean GigaRise for Free: Read Our
#    Review" paragraph right before the final (italic) paragraph.
# ------------------------------------------------------------------
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($count)

# Split the paragraph immediately preceding the last one, right before
# its own paragraph mark, so the new paragraph inherits plain (non-
# italic) formatting instead of the italic formatting of the final
# paragraph.
$prevPara = $d.Paragraphs($count - 1)
$splitPoint = $prevPara.Range.End - 1
$insPt = $d.Range($splitPoint, $splitPoint)
$insPt.InsertBefore([char]13 + "TEMP")

# Replace the just-created placeholder paragraph's contents with the
# correctly-formatted run structure via InsertXML, so we get the exact
# run layout (leading empty run + bold run) used elsewhere in the doc.
$newPara = $d.Paragraphs($count)
$newRange = $newPara.Range.Duplicate

$xmlSnippet = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
  <w:r/>
  <w:r>
    <w:rPr><w:b/></w:rPr>
    <w:t>Play Atlantean GigaRise for Free: Read Our Review</w:t>
  </w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$newRange.InsertXML($xmlSnippet)

# ------------------------------------------------------------------
# 3) Replace the text of the final paragraph (the old image-prompt
#    text) with the meta-description text, keeping its italic run
#    formatting untouched.
# ------------------------------------------------------------------
$oldImagePrompt = 'Please create a feature image for "Atlantean GigaRise" that features a happy Maya warrior with glasses in a cartoon style. The image should depict the warrior standing next to the underwater scene of Atlantis with the mountains and colonnades in the background. The warrior should be holding a trident with a smile on their face and wearing glasses. The image should be colorful and eye-catching to attract the attention of players.'
$newMetaDescription = 'Discover the features of Atlantean GigaRise, a highly volatile slot game with up to 294,912 ways to win. Play for free and read our review to learn more.'

$d.Content.Find.Execute($oldImagePrompt, $true, $false, $false, $false, $false,
                         $true, 1, $false, $newMetaDescription, 2)
